# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#   - swap the displayed country names for the two pairs of rows that
#     were re-ordered ("Santa Lucia"/"Timor Oriental" and
#     "Montserrat"/"Islas Malvinas"), including the case-count figures
#     that travelled along with "Islas Malvinas"/"Montserrat"
#   - refresh the case counters for Ucrania, Kirguistan, Afganistan and
#     Taiwan
#   - bump the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 08:48"

# --- Ucrania (row 28): totals, new cases, recovered, deaths-today, deaths ---
$ws.Range("B28").Value = 178353
$ws.Range("C28").Value = 2675
$ws.Range("E28").Value = 97258
$ws.Range("G28").Value = 26
$ws.Range("H28").Value = 3583

# --- Kirguistan (row 66): totals, new cases, active, recovered ---
$ws.Range("B66").Value = 45471
$ws.Range("C66").Value = 55
$ws.Range("D66").Value = 41682
$ws.Range("E66").Value = 2726

# --- Afganistan (row 68): totals, new cases, recovered, deaths-today, deaths ---
$ws.Range("B68").Value = 39074
$ws.Range("C68").Value = 30
$ws.Range("E68").Value = 5054
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 1444

# --- Taiwan (row 176): totals, new cases, recovered ---
$ws.Range("B176").Value = 509
$ws.Range("C176").Value = 2
$ws.Range("E176").Value = 23

# --- Rows 204/205 swap: "Timor Oriental" <-> "Santa Lucia" ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# --- Rows 214/215 swap: "Islas Malvinas" <-> "Montserrat" (with their active/deaths figures) ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
